$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3587.25
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 1449.6666
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 1449.6666
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -2697.6666
$ws.Range("H65").Value = 3587.25
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 1449.6666
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 7248.333000000001
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -13488.333
$ws.Range("H76").Value = 4114.2856
$ws.Range("I76").Value = 4560
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 4560
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -4245
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 4114.2856
$ws.Range("I79").Value = 4560
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 4560
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -3468
$ws.Range("N79").Value = -5184
$ws.Range("H87").Value = 42999.6
$ws.Range("J87").Value = 42999.6
$ws.Range("L87").Value = 42999.6
$ws.Range("N87").Value = -45495.6
$ws.Range("H90").Value = 42999.6
$ws.Range("J90").Value = 42999.6
$ws.Range("L90").Value = 128998.8
$ws.Range("N90").Value = -141478.8
$ws.Range("H98").Value = 1383
$ws.Range("I98").Value = 1456.9375
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 1456.9375
$ws.Range("L98").Value = 200
$ws.Range("M98").Value = 41.0625
$ws.Range("N98").Value = -3196
$ws.Range("H122").Value = 1383
$ws.Range("I122").Value = 1456.9375
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 4370.8125
$ws.Range("L122").Value = 600
$ws.Range("M122").Value = -1920.8125
$ws.Range("N122").Value = -5500
$ws.Range("H129").Value = 893.8039
$ws.Range("J129").Value = 899.6799999999999
$ws.Range("L129").Value = 2699.04
$ws.Range("N129").Value = -12699.04
$ws.Range("H138").Value = 2124.71
$ws.Range("I138").Value = 1036.6
$ws.Range("J138").Value = 2245.611
$ws.Range("K138").Value = 3109.8
$ws.Range("L138").Value = 6736.833
$ws.Range("M138").Value = 2030.2
$ws.Range("N138").Value = -17016.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1322.174
$ws.Range("I61").Value = 1204.8
$ws.Range("K61").Value = 1204.8
$ws.Range("M61").Value = -992.8
$ws.Range("H110").Value = 1051.3334
$ws.Range("I110").Value = 859.5185
$ws.Range("K110").Value = 859.5185
$ws.Range("M110").Value = 1185.4815
$ws.Range("H136").Value = 1322.174
$ws.Range("I136").Value = 1204.8
$ws.Range("K136").Value = 3614.4
$ws.Range("M136").Value = -1064.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3513.3635
$ws.Range("J86").Value = 6800
$ws.Range("L86").Value = 6800
$ws.Range("N86").Value = -9046
$ws.Range("H89").Value = 3513.3635
$ws.Range("J89").Value = 6800
$ws.Range("L89").Value = 34000
$ws.Range("N89").Value = -45232
$ws.Range("H99").Value = 20001132
$ws.Range("I99").Value = 23810580
$ws.Range("J99").Value = 1527.75
$ws.Range("K99").Value = 23810580
$ws.Range("L99").Value = 1527.75
$ws.Range("M99").Value = -23809082
$ws.Range("N99").Value = -4523.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 830.07465
$ws.Range("I31").Value = 731.42224
$ws.Range("J31").Value = 1031.8636
$ws.Range("K31").Value = 731.42224
$ws.Range("L31").Value = 1031.8636
$ws.Range("M31").Value = -436.42224
$ws.Range("N31").Value = -1621.8636
$ws.Range("H34").Value = 830.07465
$ws.Range("I34").Value = 731.42224
$ws.Range("J34").Value = 1031.8636
$ws.Range("K34").Value = 731.42224
$ws.Range("L34").Value = 1031.8636
$ws.Range("M34").Value = -529.42224
$ws.Range("N34").Value = -1435.8636
$ws.Range("H99").Value = 1991.6666
$ws.Range("I99").Value = 1819.6666
$ws.Range("J99").Value = 2335.6667
$ws.Range("K99").Value = 1819.6666
$ws.Range("L99").Value = 2335.6667
$ws.Range("M99").Value = -321.6666
$ws.Range("N99").Value = -5331.6667
$ws.Range("H126").Value = 1991.6666
$ws.Range("I126").Value = 1819.6666
$ws.Range("J126").Value = 2335.6667
$ws.Range("K126").Value = 5458.9998
$ws.Range("L126").Value = 7007.000100000001
$ws.Range("M126").Value = -2988.9998
$ws.Range("N126").Value = -11947.0001
$ws.Range("H132").Value = 6759.773
$ws.Range("I132").Value = 10300.272
$ws.Range("J132").Value = 3219.2727
$ws.Range("K132").Value = 30900.816
$ws.Range("L132").Value = 9657.8181
$ws.Range("M132").Value = -28370.816
$ws.Range("N132").Value = -14717.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1833.3334
$ws.Range("J59").Value = 4750
$ws.Range("L59").Value = 14250
$ws.Range("N59").Value = -15330
$ws.Range("H76").Value = 6546.5415
$ws.Range("I76").Value = 5503.25
$ws.Range("J76").Value = 6755.2
$ws.Range("K76").Value = 16509.75
$ws.Range("L76").Value = 20265.6
$ws.Range("M76").Value = -16126.75
$ws.Range("N76").Value = -21031.6
$ws.Range("H79").Value = 6546.5415
$ws.Range("I79").Value = 5503.25
$ws.Range("J79").Value = 6755.2
$ws.Range("K79").Value = 16509.75
$ws.Range("L79").Value = 20265.6
$ws.Range("M79").Value = -15183.75
$ws.Range("N79").Value = -22917.6
$ws.Range("H140").Value = 31740.297
$ws.Range("I140").Value = 47574.793
$ws.Range("J140").Value = 2507.3845
$ws.Range("K140").Value = 142724.379
$ws.Range("L140").Value = 7522.1535
$ws.Range("M140").Value = -137544.379
$ws.Range("N140").Value = -17882.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22504720
$ws.Range("I70").Value = 27782012
$ws.Range("J70").Value = 18186936
$ws.Range("K70").Value = 27782012
$ws.Range("L70").Value = 18186936
$ws.Range("M70").Value = -27781742
$ws.Range("N70").Value = -18187476
$ws.Range("H73").Value = 22504720
$ws.Range("I73").Value = 27782012
$ws.Range("J73").Value = 18186936
$ws.Range("K73").Value = 27782012
$ws.Range("L73").Value = 18186936
$ws.Range("M73").Value = -27781076
$ws.Range("N73").Value = -18188808
$ws.Range("H80").Value = 3814.6875
$ws.Range("J80").Value = 5083.7
$ws.Range("L80").Value = 5083.7
$ws.Range("N80").Value = -7079.7
$ws.Range("H83").Value = 3814.6875
$ws.Range("J83").Value = 5083.7
$ws.Range("L83").Value = 25418.5
$ws.Range("N83").Value = -35402.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H99").Value = 6380.143
$ws.Range("I99").Value = 6380.143
$ws.Range("K99").Value = 6380.143
$ws.Range("M99").Value = -4134.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 999.5714
$ws.Range("I107").Value = 799.4
$ws.Range("K107").Value = 2398.2
$ws.Range("M107").Value = -478.1999999999998
$ws.Range("H132").Value = 1772.875
$ws.Range("I132").Value = 1411.9286
$ws.Range("J132").Value = 4299.5
$ws.Range("K132").Value = 4235.7858
$ws.Range("L132").Value = 12898.5
$ws.Range("M132").Value = -1705.7858
$ws.Range("N132").Value = -17958.5
